# Adiciona o trecho final que faltava no paragrafo 1 do topico I
# ("Vem ao exame ... cargaHoraria.") -> acrescenta " objetivoProjeto "
# seguido do trecho em italico "descricaoProposta".

$d = $word.ActiveDocument

# Localiza o paragrafo que contem o texto-chave e termina em "cargaHoraria."
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -like "Vem ao exame*cargaHoraria.*") {
        $r = $p.Range
        # Colapsa para o final do paragrafo, mas antes da marca de fim de paragrafo
        $r.SetRange($r.End - 1, $r.End - 1)

        $r.InsertAfter(" objetivoProjeto ")
        $r.Collapse(0)

        $italicStart = $r.Start
        $r.InsertAfter("descricaoProposta")
        $r.Collapse(0)

        $italicRange = $d.Range($italicStart, $r.Start)
        $italicRange.Italic = 1

        break
    }
}
